$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Food
$ws.Range("B2").Value = 400
$ws.Range("C2").Value = 58.6
$ws.Range("D2").Value = 341.4

# Row 3 - Transportation
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 190

# Row 7 - Travel
$ws.Range("B7").Value = 400
$ws.Range("D7").Value = 400

# Row 9 - Medical
$ws.Range("B9").Value = 400
$ws.Range("D9").Value = 400

# Row 10 - Personal
$ws.Range("B10").Value = 400
$ws.Range("D10").Value = 400

# Row 11 - Pets
$ws.Range("C11").Value = 195.6
$ws.Range("D11").Value = 4.400000000000006
